$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $range = $ws.Range($cell)
    $range.Value = "'" + $value
    $range.ClearFormats()
}

Set-TextValue "D2" "305.87"
Set-TextValue "E2" "0.89%"
Set-TextValue "D3" "36.13"
Set-TextValue "E3" "-2.82%"
Set-TextValue "D4" "5.103"
Set-TextValue "E4" "2.15%"
Set-TextValue "D5" "0.07878"
Set-TextValue "E5" "0.80%"
Set-TextValue "D6" "2.130"
Set-TextValue "E6" "-3.22%"
Set-TextValue "D7" "7.958"
Set-TextValue "E7" "-0.69%"
Set-TextValue "D8" "0.9171"
Set-TextValue "E8" "0.13%"
Set-TextValue "D9" "0.09651"
Set-TextValue "E9" "-0.30%"
Set-TextValue "E10" "-0.69%"
Set-TextValue "D11" "0.08687"
Set-TextValue "E11" "1.78%"
Set-TextValue "D12" "0.03549"
Set-TextValue "E12" "-0.53%"
Set-TextValue "E13" "-0.37%"
Set-TextValue "D14" "0.001444"
Set-TextValue "E14" "-2.70%"
Set-TextValue "D15" "0.005672"
Set-TextValue "E15" "0.05%"
Set-TextValue "D16" "3.470"
Set-TextValue "E16" "0.33%"
Set-TextValue "D17" "4.114"
Set-TextValue "E17" "1.69%"
Set-TextValue "D18" "2.713"
Set-TextValue "E18" "14.18%"
Set-TextValue "E20" "1.26%"
Set-TextValue "D21" "5.168"
Set-TextValue "E21" "8.08%"
Set-TextValue "D23" "0.04560"
Set-TextValue "E23" "-1.22%"
Set-TextValue "D24" "0.005046"
Set-TextValue "E24" "5.51%"
Set-TextValue "D25" "0.001235"
Set-TextValue "E25" "0.12%"
Set-TextValue "E26" "14.31%"
Set-TextValue "E27" "-0.05%"
Set-TextValue "D39" "0.01854"
Set-TextValue "E39" "4.66%"
Set-TextValue "D40" "0.04777"
Set-TextValue "E40" "0.73%"
Set-TextValue "D41" "0.007454"
Set-TextValue "E41" "-7.74%"
Set-TextValue "D42" "0.1401"
Set-TextValue "E42" "0.60%"
Set-TextValue "E43" "-1.17%"
Set-TextValue "D44" "0.002231"
Set-TextValue "E44" "6.63%"
Set-TextValue "D45" "0.01105"
Set-TextValue "E45" "10.86%"
Set-TextValue "D46" "0.00006322"
Set-TextValue "E46" "3.35%"
Set-TextValue "E47" "-0.34%"
Set-TextValue "E48" "-0.02%"
Set-TextValue "D49" "47.49"
Set-TextValue "E49" "507.65%"
Set-TextValue "E50" "-25.68%"
Set-TextValue "D51" "0.00002100"
Set-TextValue "E51" "-0.34%"
